$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01251000175016376
$ws.Range("C2").Value = 0.7890050351690446
$ws.Range("D2").Value = 2.035788612216479
$ws.Range("E2").Value = 1.426810643433977
$ws.Range("F2").Value = 1.431633611568203
$ws.Range("G2").Value = 147
$ws.Range("B3").Value = 0.02111910615037019
$ws.Range("C3").Value = 0.8033144239384737
$ws.Range("D3").Value = 2.007830181931747
$ws.Range("E3").Value = 1.416979245413195
$ws.Range("F3").Value = 1.4216990518437
$ws.Range("G3").Value = 146
$ws.Range("B4").Value = 0.0285192289397949
$ws.Range("C4").Value = 0.8031463991591674
$ws.Range("D4").Value = 1.99355748821563
$ws.Range("E4").Value = 1.411933953205896
$ws.Range("F4").Value = 1.416538966206368
$ws.Range("G4").Value = 145
$ws.Range("B5").Value = 0.03163772055860532
$ws.Range("C5").Value = 0.8097221098892455
$ws.Range("D5").Value = 2.060327991514458
$ws.Range("E5").Value = 1.435384266151214
$ws.Range("F5").Value = 1.440044421046429
$ws.Range("G5").Value = 144
$ws.Range("B6").Value = 0.03800577977396699
$ws.Range("C6").Value = 0.8093017064625553
$ws.Range("D6").Value = 2.019953988652261
$ws.Range("E6").Value = 1.421250853527364
$ws.Range("F6").Value = 1.425736443153998
$ws.Range("G6").Value = 143
$ws.Range("B7").Value = 0.04592074729150665
$ws.Range("C7").Value = 0.8016490720135074
$ws.Range("D7").Value = 2.005508156182165
$ws.Range("E7").Value = 1.416159650668725
$ws.Range("F7").Value = 1.420425270060055
$ws.Range("G7").Value = 142
$ws.Range("B8").Value = 0.06000521770874007
$ws.Range("C8").Value = 0.7979248303654495
$ws.Range("D8").Value = 2.012546806653637
$ws.Range("E8").Value = 1.418642592992906
$ws.Range("F8").Value = 1.422426024816887
$ws.Range("G8").Value = 141
